# Auto-generated script to apply 2023-09-06 YTD violent crime data update
# Updates column J (year 2023 running total) across Citywide Totals,
# By Neighborhood, and 43 individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 5179
$ws.Range("J3").Value = 5512
$ws.Range("J4").Value = 1222
$ws.Range("J5").Value = 432
$ws.Range("J6").Value = 6881
$ws.Range("J7").Value = 19226

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J4").Value = 71
$ws.Range("J5").Value = 61
$ws.Range("J7").Value = 562
$ws.Range("J8").Value = 1218
$ws.Range("J10").Value = 130
$ws.Range("J11").Value = 300
$ws.Range("J12").Value = 40
$ws.Range("J14").Value = 95
$ws.Range("J19").Value = 550
$ws.Range("J20").Value = 399
$ws.Range("J21").Value = 54
$ws.Range("J22").Value = 51
$ws.Range("J23").Value = 186
$ws.Range("J29").Value = 1075
$ws.Range("J31").Value = 174
$ws.Range("J33").Value = 880
$ws.Range("J40").Value = 43
$ws.Range("J42").Value = 786
$ws.Range("J44").Value = 146
$ws.Range("J47").Value = 145
$ws.Range("J49").Value = 129
$ws.Range("J50").Value = 117
$ws.Range("J51").Value = 246
$ws.Range("J52").Value = 486
$ws.Range("J54").Value = 366
$ws.Range("J55").Value = 251
$ws.Range("J57").Value = 82
$ws.Range("J63").Value = 69
$ws.Range("J65").Value = 500
$ws.Range("J66").Value = 59
$ws.Range("J67").Value = 741
$ws.Range("J73").Value = 177
$ws.Range("J75").Value = 55
$ws.Range("J78").Value = 242
$ws.Range("J79").Value = 552
$ws.Range("J84").Value = 157
$ws.Range("J85").Value = 824
$ws.Range("J88").Value = 212
$ws.Range("J89").Value = 244
$ws.Range("J91").Value = 215
$ws.Range("J93").Value = 83
$ws.Range("J94").Value = 183
$ws.Range("J95").Value = 286
$ws.Range("J99").Value = 302
$ws.Range("J101").Value = 19226

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J3").Value = 169
$ws.Range("J4").Value = 22
$ws.Range("J6").Value = 180
$ws.Range("J7").Value = 562

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 244

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 96
$ws.Range("J6").Value = 117
$ws.Range("J7").Value = 300

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 115
$ws.Range("J6").Value = 200
$ws.Range("J7").Value = 486

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J2").Value = 35
$ws.Range("J4").Value = 7
$ws.Range("J6").Value = 32
$ws.Range("J7").Value = 95

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 339
$ws.Range("J3").Value = 371
$ws.Range("J4").Value = 73
$ws.Range("J7").Value = 1218

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 214
$ws.Range("J3").Value = 304
$ws.Range("J6").Value = 234
$ws.Range("J7").Value = 824

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 100
$ws.Range("J7").Value = 286

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J3").Value = 118
$ws.Range("J7").Value = 302

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 181
$ws.Range("J3").Value = 287
$ws.Range("J6").Value = 195
$ws.Range("J7").Value = 741

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J2").Value = 66
$ws.Range("J7").Value = 174

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J3").Value = 49
$ws.Range("J7").Value = 157

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 288
$ws.Range("J4").Value = 38
$ws.Range("J6").Value = 297
$ws.Range("J7").Value = 880

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 142
$ws.Range("J7").Value = 500

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 129

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J4").Value = 29
$ws.Range("J6").Value = 173
$ws.Range("J7").Value = 366

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 321
$ws.Range("J3").Value = 368
$ws.Range("J7").Value = 1075

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 138
$ws.Range("J6").Value = 203
$ws.Range("J7").Value = 550

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 146

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 158
$ws.Range("J6").Value = 400
$ws.Range("J7").Value = 786

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J2").Value = 67
$ws.Range("J3").Value = 78
$ws.Range("J7").Value = 242

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J4").Value = 10
$ws.Range("J6").Value = 121
$ws.Range("J7").Value = 251

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J2").Value = 51
$ws.Range("J6").Value = 47
$ws.Range("J7").Value = 186

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 215

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 54

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 155
$ws.Range("J5").Value = 15
$ws.Range("J6").Value = 153
$ws.Range("J7").Value = 552

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J6").Value = 106
$ws.Range("J7").Value = 399

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 83

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 99
$ws.Range("J7").Value = 183

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J6").Value = 67
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J3").Value = 31
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 117

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 59

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 177

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J3").Value = 61
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 212

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 61

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("J2").Value = 25
$ws.Range("J7").Value = 55

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J6").Value = 90
$ws.Range("J7").Value = 246

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J3").Value = 23
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("J2").Value = 23
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J6").Value = 25
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("J2").Value = 5
$ws.Range("J7").Value = 40
